# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The account-statement detail table (rows 16-50, columns C:G) is replaced:
# the previous records (sorted by period, then worker) are removed and a
# new set of records for the same 8 workers (sorted by worker, then by
# period descending) is written in their place. Columns B (doc type "CC")
# and H:J (always blank) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Each row: DocNumber, WorkerName, Period, ValorMora, SalarioBasico
$data = @(
    @("45524119", "DANILSA NAVARRO CUETO", "2111", 44579, 1453642),
    @("45524119", "DANILSA NAVARRO CUETO", "2110", 58146, 1453642),
    @("45524119", "DANILSA NAVARRO CUETO", "2109", 58146, 1453642),
    @("45524119", "DANILSA NAVARRO CUETO", "2108", 58146, 1453642),
    @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2111", 44579, 1453642),
    @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2110", 58146, 1453642),
    @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2109", 58146, 1453642),
    @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2108", 58146, 1453642),
    @("9282469", "RUBEN SOTO MARTINEZ", "2111", 27861, 908526),
    @("9282469", "RUBEN SOTO MARTINEZ", "2110", 36341, 908526),
    @("9282469", "RUBEN SOTO MARTINEZ", "2109", 36341, 908526),
    @("9282469", "RUBEN SOTO MARTINEZ", "2108", 36341, 908526),
    @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2111", 44579, 1453642),
    @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2110", 58146, 1453642),
    @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2109", 58146, 1453642),
    @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2108", 42640, 1453642),
    @("73183791", "YESID QUINTANA TORRES", "2111", 27861, 908526),
    @("73183791", "YESID QUINTANA TORRES", "2110", 36341, 908526),
    @("73183791", "YESID QUINTANA TORRES", "2109", 36341, 908526),
    @("73183791", "YESID QUINTANA TORRES", "2108", 36341, 908526),
    @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2111", 27861, 908526),
    @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2110", 36341, 908526),
    @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2109", 36341, 908526),
    @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2108", 36341, 908526),
    @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2107", 36341, 908526),
    @("1193561465", "HERNAN DARIO VILA NORIEGA", "2111", 55723, 1817052),
    @("1193561465", "HERNAN DARIO VILA NORIEGA", "2110", 72682, 1817052),
    @("1193561465", "HERNAN DARIO VILA NORIEGA", "2109", 72682, 1817052),
    @("1193561465", "HERNAN DARIO VILA NORIEGA", "2108", 72682, 1817052),
    @("1193561465", "HERNAN DARIO VILA NORIEGA", "2107", 72682, 1817052),
    @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2111", 62688, 2044184),
    @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2110", 81767, 2044184),
    @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2109", 81767, 2044184),
    @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2108", 81767, 2044184),
    @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2107", 81767, 2044184)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 3).Value = $row[0]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[2]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[3]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[4]   # G: Salario Basico
}
